# Updates the "Permissões" sheet: refreshes the function/procedure names
# in the right-hand "Permissões de Execução" tables to match the new
# naming convention (sys_call.*_t_<table>_r_<role>()) and renames the
# header column from "Nome da função" to "Nome da função/procedure".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Permissões")

# --- Table 1 (auth.users) : rows 13-18 --------------------------------
$ws.Range("M13").Value = "Nome da função/procedure"

$ws.Range("M14").Value = "sys_call.view_insert_function_t_call_user_r_sys_usr()"
$ws.Range("N14").Value = "usr"
$ws.Range("O14").Value = "sys_usr"

$ws.Range("M15").Value = "sys_call.function_login_t_call_user_r_sys_usr()"
$ws.Range("N15").Value = "usr"
$ws.Range("O15").Value = "sys_usr"

$ws.Range("M16").Value = "sys_call.function_select_all_t_call_user_r_sys_reg_usr()"
$ws.Range("N16").Value = "reg_usr"
$ws.Range("O16").Value = "sys_reg_usr"

$ws.Range("M17").Value = "sys_call.view_update_function_t_call_user_r_sys_reg_usr()"
$ws.Range("N17").Value = "reg_usr"
$ws.Range("O17").Value = "sys_reg_usr"

$ws.Range("M18").Value = "sys_call.procedure_delete_t_call_user_r_sys_reg_usr()"
$ws.Range("N18").Value = "reg_usr"
$ws.Range("O18").Value = "sys_reg_usr"

# --- Table 2 (sys_call.contact) : rows 28-32 ---------------------------
$ws.Range("M28").Value = "Nome da função/procedure"

$ws.Range("M29").Value = "sys_call.view_insert_function_t_contact_r_sys_reg_usr()"
$ws.Range("M30").Value = "sys_call.function_select_t_contact_r_sys_reg_usr()"
$ws.Range("M31").Value = "sys_call.view_update_function_t_contact_r_sys_reg_usr()"
$ws.Range("M32").Value = "sys_call.procedure_delete_t_contact_r_sys_reg_usr()"

# --- Table 3 (public.call_history) : header only -----------------------
$ws.Range("M40").Value = "Nome da função/procedure"

$ws.Columns("M").ColumnWidth = 54.1
